# Levine plot for report:
# Insert a new "PA" column right after "Age" (new column C), duplicating the
# existing "ph_age" column values, shifting Albumin..acceleration one column
# to the right (D..O). Matches the author's commit "levine plot for report".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C (Albumin). This shifts
# all subsequent columns (Albumin .. acceleration) one place to the right
# and keeps the per-row number formatting (style) that column B already had,
# since Excel's column insert inherits formatting from the column to the left.
$ws.Range("C1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value2 = "PA"

# Fill the new "PA" column with the same values as "ph_age" (now column N,
# since it shifted right by one place along with the other measurement
# columns).
for ($r = 2; $r -le 34; $r++) {
    $phAge = $ws.Cells.Item($r, 14).Value2
    $ws.Cells.Item($r, 3).Value2 = $phAge
}

# Restore the tab selection to the cell the author left active.
[void]$ws.Range("C40").Select()
